# Daily attendance processing - 2025-12-16 11:00:38
#
# The "Recorded By" column (G) on the session-analysis sheet lists the
# accounts that touched each attendance record, comma-separated. This pass
# re-normalizes the ordering of that list for the affected rows by swapping
# the last two comma-separated entries (e.g. "dnasr281@gmail.com, System"
# -> "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows in column G whose "Recorded By" list needs the last two entries swapped.
$rowsToFix = @(
    2, 3, 6, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26, 28, 29, 32,
    36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 54, 55, 58, 62, 63,
    64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 87, 90, 92,
    93, 94, 96, 99, 101, 109, 110, 111, 112, 113, 116, 118, 119, 120, 122, 125,
    127, 135, 136, 137, 138, 139, 142, 144, 145, 146, 148, 151, 153
)

foreach ($row in $rowsToFix) {
    $cell = $ws.Cells.Item($row, 7)
    $parts = $cell.Value2.Split(",") | ForEach-Object { $_.Trim() }
    $lastIdx = $parts.Count - 1
    $secondLastIdx = $parts.Count - 2

    $tmp = $parts[$lastIdx]
    $parts[$lastIdx] = $parts[$secondLastIdx]
    $parts[$secondLastIdx] = $tmp

    $cell.Value = [string]::Join(", ", $parts)
}
